# Update "想去人数" (want-to-go count) figures in the 展览 and 全部类型
# sheets, as refreshed by the gh-pages data generation job.

$wb = $excel.ActiveWorkbook

# Row (by sheet r="N") -> new F-column value
$updates = @{
    3  = 1074
    8  = 11136
    9  = 4273
    11 = 20
    12 = 14
    13 = 2496
    15 = 90
    17 = 155
    19 = 11213
    20 = 11060
    22 = 36
}

foreach ($sheetName in @("展览", "全部类型")) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($row in $updates.Keys) {
        $ws.Range("F$row").Value = $updates[$row]
    }
}
